# Vehicle Loan Default Prediction.pptx — insert 5 new blank "Title and
# Content" slides at the positions required by the target deck, leaving
# the existing 5 slides (Title/Team, Problem Definition, Suggested
# Solution and EDA, Algorithms/Solution/Conclusions, Thanks!) with their
# original content untouched, just shifted to new slide positions.
#
# Before: 1 Title/Team | 2 Problem Definition | 3 Suggested Solution |
#         4 Algorithms  | 5 Thanks!
# After:  1 Title/Team | 2 Problem Definition | 3 NEW | 4 NEW |
#         5 Suggested Solution | 6 NEW | 7 Algorithms | 8 NEW | 9 NEW |
#         10 Thanks!
#
# ppLayoutText (layout index 2 in the slide master's layout list) is the
# "Title and Content" layout, giving each new slide a Title placeholder
# plus a body/content placeholder, both left empty.

$p = $ppt.ActivePresentation

$p.Slides.Add(3, 2)  | Out-Null
$p.Slides.Add(4, 2)  | Out-Null
$p.Slides.Add(6, 2)  | Out-Null
$p.Slides.Add(8, 2)  | Out-Null
$p.Slides.Add(9, 2)  | Out-Null

Write-Output "Slide count: $($p.Slides.Count)"
